$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $value) {
    # Force the cell to be treated as text so Excel does not silently
    # convert price-looking strings (e.g. "312.71") into numeric values
    # (which would lose trailing zeros / introduce floating point noise).
    $c = $ws.Cells.Item($row, 4)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Price 2 "28.289.47"
Set-Volume 2 "-0.30%"

# Row 3 - Ethereum
Set-Price 3 "1.810.56"
Set-Volume 3 "-0.55%"

# Row 4 - TetherUSD
Set-Volume 4 "-0.16%"

# Row 5 - BNB
Set-Price 5 "312.71"
Set-Volume 5 "-0.78%"

# Row 6 - USDC
Set-Volume 6 "-0.14%"

# Row 7 - XRP
Set-Price 7 "0.5134"
Set-Volume 7 "-2.18%"

# Row 8 - Cardano
Set-Price 8 "0.3933"
Set-Volume 8 "+2.24%"

# Row 9 - Dogecoin
Set-Price 9 "0.07821"
Set-Volume 9 "-2.83%"

# Row 10 - Polygon
Set-Price 10 "1.109"
Set-Volume 10 "-0.30%"

# Row 11 - OKB
Set-Volume 11 "-1.97%"

# Row 12 - Polkadot
Set-Price 12 "6.369"
Set-Volume 12 "-0.50%"

# Row 13 - BinanceUSD
Set-Volume 13 "-0.13%"

# Row 14 - Solana
Set-Price 14 "20.47"
Set-Volume 14 "-1.79%"

# Row 15 - Chainlink
Set-Price 15 "7.338"
Set-Volume 15 "-1.08%"

# Row 16 - WrappedEther
Set-Price 16 "1.801.24"
Set-Volume 16 "-0.59%"

# Row 17 - Litecoin
Set-Price 17 "92.76"
Set-Volume 17 "-1.77%"

# Row 18 - ShibaInu
Set-Price 18 "0.00001079"
Set-Volume 18 "-2.14%"

# Row 19 - TRON
Set-Price 19 "0.06583"
Set-Volume 19 "-0.85%"

# Row 20 - Dai
Set-Price 20 "1.001"
Set-Volume 20 "-0.14%"

# Row 21 - Avalanche
Set-Price 21 "17.33"
Set-Volume 21 "-1.59%"

# Row 22 - Uniswap
Set-Price 22 "6.014"
Set-Volume 22 "-0.05%"

# Row 23 - WrappedBTC
Set-Price 23 "28.332.07"

# Row 24 - Cosmos
Set-Volume 24 "-1.88%"

# Row 25 - Toncoin
Set-Price 25 "2.230"
Set-Volume 25 "-0.79%"

# Row 26 - Monero
Set-Price 26 "160.19"
Set-Volume 26 "+0.77%"

# Row 27 - LidoDAOToken
Set-Price 27 "2.463"
Set-Volume 27 "+2.33%"

# Row 28 - EthereumClassic
Set-Volume 28 "-1.74%"

# Row 29 - WrappedliquidstakedEther2.0
Set-Price 29 "2.016.17"
Set-Volume 29 "-0.47%"

# Row 30 - BitcoinCash
Set-Price 30 "127.47"
Set-Volume 30 "+2.35%"

# Row 31 - Stellar
Set-Price 31 "0.1095"
Set-Volume 31 "-1.54%"

# Row 32 - ImmutableX
Set-Price 32 "1.061"
Set-Volume 32 "-1.52%"

# Row 33 - HuobiToken
Set-Price 33 "3.654"
Set-Volume 33 "-0.69%"

# Row 34 - Filecoin
Set-Price 34 "5.577"
Set-Volume 34 "-1.56%"

# Row 35 - Hedera
Set-Price 35 "0.07139"
Set-Volume 35 "-2.66%"

# Row 36 - FraxShare
Set-Price 36 "9.132"
Set-Volume 36 "+4.80%"

# Row 37 - VeChain
Set-Price 37 "0.02352"
Set-Volume 37 "+0.42%"

# Row 38 - Algorand
Set-Price 38 "0.2177"
Set-Volume 38 "-1.09%"

# Row 39 - Aptos
Set-Price 39 "11.61"
Set-Volume 39 "-4.66%"

# Row 40 - InternetComputer(DFINITY)
Set-Price 40 "5.016"
Set-Volume 40 "-2.04%"

# Row 41 - TheSandbox
Set-Price 41 "0.6182"
Set-Volume 41 "-1.92%"

# Row 42 - Frax
Set-Price 42 "1.001"
Set-Volume 42 "-0.17%"

# Row 43 - TrustWalletToken
Set-Price 43 "1.159"
Set-Volume 43 "-1.76%"

# Row 44 - EnergySwap
Set-Price 44 "13.22"
Set-Volume 44 "-1.04%"

# Row 45 - Decentraland
Set-Price 45 "0.5965"
Set-Volume 45 "-2.59%"

# Row 46 - WEMIXTOKEN
Set-Price 46 "1.304"
Set-Volume 46 "-5.75%"

# Row 47 - PancakeSwap
Set-Volume 47 "-1.58%"

# Row 48 - Quant
Set-Price 48 "125.36"
Set-Volume 48 "-1.31%"

# Row 49 - EOS
Set-Price 49 "1.211"
Set-Volume 49 "-0.08%"

# Row 50 - NEARProtocol
Set-Volume 50 "-2.24%"

# Row 51 - Cronos
Set-Volume 51 "-1.38%"
